$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting the existing D:K data to F:M
$ws.Columns("D:E").Insert()

# Carry over the number formatting (date format for header rows, #,##0 for
# data rows) from the columns that used to be D:E (now shifted to F:G) into
# the freshly inserted D:E columns.
$ws.Columns("F:G").Copy()
$ws.Columns("D:E").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("F7").Value = 43281
$ws.Range("D8").Value = "NA"
$ws.Range("E8").Value = 9700
$ws.Range("F8").Value = 18200
$ws.Range("D9").Value = 6200
$ws.Range("E9").Value = 4800
$ws.Range("F9").Value = 9800
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = 4900
$ws.Range("F10").Value = 8400
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("D17").Value = -11700
$ws.Range("E17").Value = 8300
$ws.Range("F17").Value = 16700
$ws.Range("D18").Value = "NA"
$ws.Range("E18").Value = 1400
$ws.Range("F18").Value = 1500
$ws.Range("D20").Value = "NA"
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = 1600
$ws.Range("F21").Value = 2000
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 100
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 1400
$ws.Range("F23").Value = 1400
$ws.Range("D24").Value = -100
$ws.Range("E24").Value = 300
$ws.Range("F24").Value = 400
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("D26").Value = "NA"
$ws.Range("E26").Value = 1100
$ws.Range("F26").Value = 1100
$ws.Range("D27").Value = "NA"
$ws.Range("E27").Value = 1100
$ws.Range("F27").Value = 1100
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("D32").Value = "NA"
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("D33").Value = "NA"
$ws.Range("E33").Value = 1100
$ws.Range("F33").Value = 1100
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("D35").Value = "NA"
$ws.Range("E35").Value = 1100
$ws.Range("F35").Value = 1100
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("F38").Value = 43281
$ws.Range("D41").Value = 6200
$ws.Range("E41").Value = 5600
$ws.Range("F41").Value = 5000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("D43").Value = 3600
$ws.Range("E43").Value = 3500
$ws.Range("F43").Value = 4200
$ws.Range("D44").Value = 7800
$ws.Range("E44").Value = 8400
$ws.Range("F44").Value = 7700
$ws.Range("D45").Value = 500
$ws.Range("E45").Value = 300
$ws.Range("F45").Value = 400
$ws.Range("D46").Value = 18200
$ws.Range("E46").Value = 17900
$ws.Range("F46").Value = 17300
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("D48").Value = 1000
$ws.Range("E48").Value = 1000
$ws.Range("F48").Value = 900
$ws.Range("D49").Value = 7000
$ws.Range("E49").Value = 7200
$ws.Range("F49").Value = 7400
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("D52").Value = 300
$ws.Range("E52").Value = 500
$ws.Range("F52").Value = 500
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("D54").Value = 26500
$ws.Range("E54").Value = 26600
$ws.Range("F54").Value = 26100
$ws.Range("D57").Value = 1800
$ws.Range("E57").Value = 1900
$ws.Range("F57").Value = 2500
$ws.Range("D58").Value = "NA"
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 0
$ws.Range("D59").Value = 600
$ws.Range("E59").Value = 800
$ws.Range("F59").Value = 1100
$ws.Range("D60").Value = 2400
$ws.Range("E60").Value = 2600
$ws.Range("F60").Value = 3600
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("F61").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("F62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("D66").Value = 2400
$ws.Range("E66").Value = 2600
$ws.Range("F66").Value = 3600
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("D72").Value = 15800
$ws.Range("E72").Value = 15700
$ws.Range("F72").Value = 14600
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("D76").Value = 24100
$ws.Range("E76").Value = 24000
$ws.Range("F76").Value = 22500
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("F80").Value = 43281
$ws.Range("D81").Value = "NA"
$ws.Range("E81").Value = 1100
$ws.Range("F81").Value = 1100
$ws.Range("D83").Value = 200
$ws.Range("E83").Value = 200
$ws.Range("F83").Value = 400
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("D89").Value = 600
$ws.Range("E89").Value = 400
$ws.Range("F89").Value = 2000
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = -200
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = "NA"
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("D94").Value = 0
$ws.Range("E94").Value = -200
$ws.Range("F94").Value = 0
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("F96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("D100").Value = 0
$ws.Range("E100").Value = 400
$ws.Range("F100").Value = -1100
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 0
$ws.Range("D102").Value = 600
$ws.Range("E102").Value = 600
$ws.Range("F102").Value = 900
